# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml ("Office Theme" colours)  <->  ppt/theme/theme2.xml ("Integral" colours)
# Everything else (fonts, format scheme, relationships) is untouched - only the
# <a:clrScheme> RGB values travel between the two parts.
#
# The slide master / the whole presentation is wired (via the package rels) to
# ppt/theme/theme2.xml, so that is the part whose colours are actually seen
# throughout the deck. We reproduce the effect of the swap on the object model
# by pushing the "Office Theme" palette into that live theme's colour scheme,
# the same way PowerPoint itself rewrites the clrScheme in place when a user
# picks a different theme from the Design gallery.

function Get-RgbValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches MsoThemeColorSchemeIndex:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000",  # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $tcs.Item($i).RGB = Get-RgbValue $officeThemeHex[$i - 1]
}
